$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = -1
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "Datos incorrectos"
$ws.Range("E10").Value = "Datos incorrectos"
$ws.Range("F10").Value = "Pasado"
